$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")
$ws.Activate()
$ws.Rows(2).Delete()
$ws.Range("L6").Select()
